$wb = $excel.ActiveWorkbook
$wsResumen = $wb.Worksheets.Item("Resumen")
$wsSolucion = $wb.Worksheets.Item("Solucion")
$wsMetricas = $wb.Worksheets.Item("Metricas")

$wsResumen.Range("B2").Value = "Z3"
$wsResumen.Range("C2").Value = 543.0867499189102
$wsSolucion.Range("A2").Value = "Pedido_77"
$wsSolucion.Range("A3").Value = "Pedido_73"
$wsSolucion.Range("A4").Value = "Pedido_37"
$wsSolucion.Range("A5").Value = "Pedido_31"
$wsSolucion.Range("A6").Value = "Pedido_29"
$wsSolucion.Range("A7").Value = "Pedido_24"
$wsSolucion.Range("B7").Value = "S011"
$wsSolucion.Range("A8").Value = "Pedido_80"
$wsSolucion.Range("B8").Value = "S071"
$wsSolucion.Range("A9").Value = "Pedido_60"
$wsSolucion.Range("B9").Value = "S051"
$wsSolucion.Range("A10").Value = "Pedido_76"
$wsSolucion.Range("B10").Value = "S002"
$wsSolucion.Range("A11").Value = "Pedido_11"
$wsSolucion.Range("B11").Value = "S022"
$wsSolucion.Range("A12").Value = "Pedido_15"
$wsSolucion.Range("B12").Value = "S062"
$wsSolucion.Range("A13").Value = "Pedido_6"
$wsSolucion.Range("B13").Value = "S042"
$wsSolucion.Range("A14").Value = "Pedido_20"
$wsSolucion.Range("B14").Value = "S012"
$wsSolucion.Range("A15").Value = "Pedido_17"
$wsSolucion.Range("B15").Value = "S072"
$wsSolucion.Range("A16").Value = "Pedido_61"
$wsSolucion.Range("B16").Value = "S052"
$wsSolucion.Range("A17").Value = "Pedido_78"
$wsSolucion.Range("A18").Value = "Pedido_26"
$wsSolucion.Range("B19").Value = "S023"
$wsSolucion.Range("A20").Value = "Pedido_42"
$wsSolucion.Range("A21").Value = "Pedido_3"
$wsSolucion.Range("B21").Value = "S003"
$wsSolucion.Range("A22").Value = "Pedido_36"
$wsSolucion.Range("B22").Value = "S073"
$wsSolucion.Range("A23").Value = "Pedido_39"
$wsSolucion.Range("B23").Value = "S033"
$wsSolucion.Range("A24").Value = "Pedido_5"
$wsSolucion.Range("B24").Value = "S053"
$wsSolucion.Range("A25").Value = "Pedido_64"
$wsSolucion.Range("B25").Value = "S064"
$wsSolucion.Range("A26").Value = "Pedido_21"
$wsSolucion.Range("B26").Value = "S013"
$wsSolucion.Range("A27").Value = "Pedido_46"
$wsSolucion.Range("B27").Value = "S024"
$wsSolucion.Range("A28").Value = "Pedido_33"
$wsSolucion.Range("B28").Value = "S044"
$wsSolucion.Range("A29").Value = "Pedido_54"
$wsSolucion.Range("A30").Value = "Pedido_44"
$wsSolucion.Range("B30").Value = "S074"
$wsSolucion.Range("A31").Value = "Pedido_40"
$wsSolucion.Range("A32").Value = "Pedido_10"
$wsSolucion.Range("A33").Value = "Pedido_19"
$wsSolucion.Range("B33").Value = "S065"
$wsSolucion.Range("A34").Value = "Pedido_56"
$wsSolucion.Range("B34").Value = "S025"
$wsSolucion.Range("A35").Value = "Pedido_9"
$wsSolucion.Range("B35").Value = "S054"
$wsSolucion.Range("A36").Value = "Pedido_30"
$wsSolucion.Range("B36").Value = "S075"
$wsSolucion.Range("A37").Value = "Pedido_38"
$wsSolucion.Range("A38").Value = "Pedido_22"
$wsSolucion.Range("B38").Value = "S035"
$wsSolucion.Range("A39").Value = "Pedido_79"
$wsSolucion.Range("B39").Value = "S045"
$wsSolucion.Range("A40").Value = "Pedido_52"
$wsSolucion.Range("B40").Value = "S015"
$wsSolucion.Range("A41").Value = "Pedido_1"
$wsSolucion.Range("B41").Value = "S066"
$wsSolucion.Range("B42").Value = "S055"
$wsSolucion.Range("A43").Value = "Pedido_53"
$wsSolucion.Range("B43").Value = "S026"
$wsSolucion.Range("A44").Value = "Pedido_32"
$wsSolucion.Range("B44").Value = "S076"
$wsSolucion.Range("A45").Value = "Pedido_71"
$wsSolucion.Range("B45").Value = "S046"
$wsSolucion.Range("B46").Value = "S006"
$wsSolucion.Range("A47").Value = "Pedido_66"
$wsSolucion.Range("B47").Value = "S036"
$wsSolucion.Range("A48").Value = "Pedido_18"
$wsSolucion.Range("B48").Value = "S067"
$wsSolucion.Range("A49").Value = "Pedido_69"
$wsSolucion.Range("A50").Value = "Pedido_28"
$wsSolucion.Range("B50").Value = "S027"
$wsSolucion.Range("A51").Value = "Pedido_70"
$wsSolucion.Range("A52").Value = "Pedido_7"
$wsSolucion.Range("B52").Value = "S047"
$wsSolucion.Range("A53").Value = "Pedido_43"
$wsSolucion.Range("B53").Value = "S077"
$wsSolucion.Range("A54").Value = "Pedido_75"
$wsSolucion.Range("B54").Value = "S007"
$wsSolucion.Range("A55").Value = "Pedido_55"
$wsSolucion.Range("B55").Value = "S037"
$wsSolucion.Range("A56").Value = "Pedido_51"
$wsSolucion.Range("B56").Value = "S017"
$wsSolucion.Range("A57").Value = "Pedido_47"
$wsSolucion.Range("B57").Value = "S057"
$wsSolucion.Range("A58").Value = "Pedido_59"
$wsSolucion.Range("B58").Value = "S068"
$wsSolucion.Range("A59").Value = "Pedido_34"
$wsSolucion.Range("A60").Value = "Pedido_65"
$wsSolucion.Range("B60").Value = "S008"
$wsSolucion.Range("A61").Value = "Pedido_12"
$wsSolucion.Range("B61").Value = "S048"
$wsSolucion.Range("A62").Value = "Pedido_16"
$wsSolucion.Range("B62").Value = "S018"
$wsSolucion.Range("A63").Value = "Pedido_48"
$wsSolucion.Range("B63").Value = "S078"
$wsSolucion.Range("A64").Value = "Pedido_25"
$wsSolucion.Range("B64").Value = "S009"
$wsSolucion.Range("A65").Value = "Pedido_68"
$wsSolucion.Range("B65").Value = "S038"
$wsSolucion.Range("A66").Value = "Pedido_8"
$wsSolucion.Range("B66").Value = "S058"
$wsSolucion.Range("A67").Value = "Pedido_23"
$wsSolucion.Range("A68").Value = "Pedido_4"
$wsSolucion.Range("B68").Value = "S069"
$wsSolucion.Range("A69").Value = "Pedido_74"
$wsSolucion.Range("B69").Value = "S019"
$wsSolucion.Range("A70").Value = "Pedido_62"
$wsSolucion.Range("B70").Value = "S049"
$wsSolucion.Range("A71").Value = "Pedido_58"
$wsSolucion.Range("B71").Value = "S039"
$wsSolucion.Range("A72").Value = "Pedido_49"
$wsSolucion.Range("B72").Value = "S079"
$wsSolucion.Range("A73").Value = "Pedido_67"
$wsSolucion.Range("A74").Value = "Pedido_13"
$wsSolucion.Range("B74").Value = "S010"
$wsSolucion.Range("A75").Value = "Pedido_2"
$wsSolucion.Range("B75").Value = "S030"
$wsSolucion.Range("A76").Value = "Pedido_63"
$wsSolucion.Range("B76").Value = "S070"
$wsSolucion.Range("A77").Value = "Pedido_57"
$wsSolucion.Range("B77").Value = "S050"
$wsSolucion.Range("A78").Value = "Pedido_45"
$wsSolucion.Range("B78").Value = "S060"
$wsSolucion.Range("A79").Value = "Pedido_72"
$wsSolucion.Range("B79").Value = "S080"
$wsSolucion.Range("A80").Value = "Pedido_41"
$wsSolucion.Range("B80").Value = "S020"
$wsSolucion.Range("A81").Value = "Pedido_27"
$wsSolucion.Range("B81").Value = "S040"
$wsMetricas.Range("B2").Value = 539.5760514650233
$wsMetricas.Range("B3").Value = 541.5462482430532
$wsMetricas.Range("B4").Value = 543.0867499189102
$wsMetricas.Range("B5").Value = 515.6891555843875

Write-Output "done"
